# Auto-generated Excel COM-interop script
# Applies profit-recalculation updates to H:N columns across 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H70").Value = 2641.111
$ws.Range("I70").Value = 2728.3333
$ws.Range("J70").Value = 2466.6667
$ws.Range("K70").Value = 8184.999899999999
$ws.Range("L70").Value = 7400.000100000001
$ws.Range("M70").Value = -7914.999899999999
$ws.Range("N70").Value = -7940.000100000001
$ws.Range("H73").Value = 2641.111
$ws.Range("I73").Value = 2728.3333
$ws.Range("J73").Value = 2466.6667
$ws.Range("K73").Value = 8184.999899999999
$ws.Range("L73").Value = 7400.000100000001
$ws.Range("M73").Value = -7248.999899999999
$ws.Range("N73").Value = -9272.000100000001
$ws.Range("H107").Value = 86178.89
$ws.Range("J107").Value = 192248.75
$ws.Range("L107").Value = 192248.75
$ws.Range("N107").Value = -196088.75
$ws.Range("H121").Value = 19999
$ws.Range("J121").Value = 19999
$ws.Range("L121").Value = 59997
$ws.Range("N121").Value = -63491
$ws.Range("H129").Value = 1958.3158
$ws.Range("J129").Value = 3970.7144
$ws.Range("L129").Value = 11912.1432
$ws.Range("N129").Value = -21912.1432

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17644.61
$ws.Range("I32").Value = 3013.791
$ws.Range("J32").Value = 507777
$ws.Range("K32").Value = 3013.791
$ws.Range("L32").Value = 507777
$ws.Range("M32").Value = -2726.791
$ws.Range("N32").Value = -508351
$ws.Range("H97").Value = 780.6
$ws.Range("I97").Value = 840.6667
$ws.Range("K97").Value = 840.6667
$ws.Range("M97").Value = -344.6667
$ws.Range("H122").Value = 2904.8333
$ws.Range("I122").Value = 2995.0588
$ws.Range("J122").Value = 2685.7144
$ws.Range("K122").Value = 8985.1764
$ws.Range("L122").Value = 8057.1432
$ws.Range("M122").Value = -6535.1764
$ws.Range("N122").Value = -12957.1432
$ws.Range("H132").Value = 1915.5349
$ws.Range("I132").Value = 1133.1714
$ws.Range("K132").Value = 3399.5142
$ws.Range("M132").Value = -869.5141999999996
$ws.Range("H140").Value = 70900.375
$ws.Range("J140").Value = 70900.375
$ws.Range("L140").Value = 70900.375
$ws.Range("N140").Value = -81260.375
$ws.Range("H141").Value = 144250
$ws.Range("J141").Value = 144250
$ws.Range("L141").Value = 144250
$ws.Range("N141").Value = -154610

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620
$ws.Range("H107").Value = 64496.938
$ws.Range("J107").Value = 2278.4
$ws.Range("L107").Value = 2278.4
$ws.Range("N107").Value = -6118.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3421.0908
$ws.Range("I132").Value = 3398.611
$ws.Range("K132").Value = 10195.833
$ws.Range("M132").Value = -7665.832999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("K63").Value = 6000
$ws.Range("M63").Value = -5251
$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("K66").Value = 18000
$ws.Range("M66").Value = -14256
$ws.Range("H68").Value = 1643.1428
$ws.Range("I68").Value = 1739.8
$ws.Range("J68").Value = 1401.5
$ws.Range("K68").Value = 5219.4
$ws.Range("L68").Value = 4204.5
$ws.Range("M68").Value = -4408.4
$ws.Range("N68").Value = -5826.5
$ws.Range("H70").Value = 4097
$ws.Range("I70").Value = 1661.6666
$ws.Range("J70").Value = 7750
$ws.Range("K70").Value = 4984.9998
$ws.Range("L70").Value = 23250
$ws.Range("M70").Value = -4669.9998
$ws.Range("N70").Value = -23880
$ws.Range("H71").Value = 1643.1428
$ws.Range("I71").Value = 1739.8
$ws.Range("J71").Value = 1401.5
$ws.Range("K71").Value = 15658.2
$ws.Range("L71").Value = 12613.5
$ws.Range("M71").Value = -11602.2
$ws.Range("N71").Value = -20725.5
$ws.Range("H73").Value = 4097
$ws.Range("I73").Value = 1661.6666
$ws.Range("J73").Value = 7750
$ws.Range("K73").Value = 4984.9998
$ws.Range("L73").Value = 23250
$ws.Range("M73").Value = -3892.9998
$ws.Range("N73").Value = -25434
$ws.Range("H80").Value = 2025
$ws.Range("J80").Value = 2033.3334
$ws.Range("L80").Value = 6100.0002
$ws.Range("N80").Value = -7972.0002
$ws.Range("H83").Value = 2025
$ws.Range("J83").Value = 2033.3334
$ws.Range("L83").Value = 18300.0006
$ws.Range("N83").Value = -27660.0006
$ws.Range("H113").Value = 1383.3636
$ws.Range("J113").Value = 1641.7142
$ws.Range("L113").Value = 4925.142599999999
$ws.Range("N113").Value = -9265.142599999999
$ws.Range("H129").Value = 2571.8333
$ws.Range("I129").Value = 1056
$ws.Range("J129").Value = 3154.8462
$ws.Range("K129").Value = 3168
$ws.Range("L129").Value = 9464.5386
$ws.Range("M129").Value = 1832
$ws.Range("N129").Value = -19464.5386

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6073.4707
$ws.Range("I70").Value = 6110.5454
$ws.Range("J70").Value = 6005.5
$ws.Range("K70").Value = 6110.5454
$ws.Range("L70").Value = 6005.5
$ws.Range("M70").Value = -5840.5454
$ws.Range("N70").Value = -6545.5
$ws.Range("H73").Value = 6073.4707
$ws.Range("I73").Value = 6110.5454
$ws.Range("J73").Value = 6005.5
$ws.Range("K73").Value = 6110.5454
$ws.Range("L73").Value = 6005.5
$ws.Range("M73").Value = -5174.5454
$ws.Range("N73").Value = -7877.5
$ws.Range("H102").Value = 2283.476
$ws.Range("I102").Value = 2107.7222
$ws.Range("K102").Value = 2107.7222
$ws.Range("M102").Value = -485.7222000000002
$ws.Range("H126").Value = 4499.5
$ws.Range("I126").Value = 4531.6
$ws.Range("J126").Value = 4339
$ws.Range("K126").Value = 13594.8
$ws.Range("L126").Value = 13017
$ws.Range("M126").Value = -11124.8
$ws.Range("N126").Value = -17957

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15230.73
$ws.Range("I7").Value = 55199.8
$ws.Range("K7").Value = 55199.8
$ws.Range("M7").Value = -55087.8
$ws.Range("H40").Value = 7999.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 7999.5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 7999.5
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -8271.5
$ws.Range("H46").Value = 41297.363
$ws.Range("I46").Value = 142124
$ws.Range("J46").Value = 3487.375
$ws.Range("K46").Value = 142124
$ws.Range("L46").Value = 3487.375
$ws.Range("M46").Value = -141936
$ws.Range("N46").Value = -3863.375
$ws.Range("H122").Value = 5110.622
$ws.Range("I122").Value = 5305.794
$ws.Range("J122").Value = 4507.364
$ws.Range("K122").Value = 15917.382
$ws.Range("L122").Value = 13522.092
$ws.Range("M122").Value = -13467.382
$ws.Range("N122").Value = -18422.092
$ws.Range("H126").Value = 15230.73
$ws.Range("I126").Value = 55199.8
$ws.Range("K126").Value = 165599.4
$ws.Range("M126").Value = -163129.4
$ws.Range("H132").Value = 3824.1755
$ws.Range("I132").Value = 3214.4524
$ws.Range("K132").Value = 9643.3572
$ws.Range("M132").Value = -7113.3572

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 12399.728
$ws.Range("I51").Value = 1412.375
$ws.Range("J51").Value = 18678.215
$ws.Range("K51").Value = 1412.375
$ws.Range("L51").Value = 18678.215
$ws.Range("M51").Value = -902.375
$ws.Range("N51").Value = -19698.215
$ws.Range("H100").Value = 2981.1875
$ws.Range("I100").Value = 3021.3572
$ws.Range("J100").Value = 2700
$ws.Range("K100").Value = 6042.7144
$ws.Range("L100").Value = 5400
$ws.Range("M100").Value = -5501.7144
$ws.Range("N100").Value = -6482
